# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the "bec8927b..."
# and "cce2497c..." files are now ready for handoff (instead of "handed
# back"), records the new handoff timestamps, and records an error detail
# explaining that the handback file version is stale.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

$errBec = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66e0868eb567f83d8680587cdd012948ce9741b7/e2e/bec8927b-3d2b-464a-9142-5459724d49a8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93216558d6bd9c8e183301c9c1042ee7ac93a0de/e2e/bec8927b-3d2b-464a-9142-5459724d49a8.md."
$errCce = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66e0868eb567f83d8680587cdd012948ce9741b7/e2e/cce2497c-83a3-4551-b61d-685b84079eb2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93216558d6bd9c8e183301c9c1042ee7ac93a0de/e2e/cce2497c-83a3-4551-b61d-685b84079eb2.md."

# The engine round-trips ColumnWidth (character units) through a pixel-based
# conversion that adds a flat 5/6 when writing the stored OOXML "width"
# attribute. Subtracting that offset here makes the persisted width land on
# exactly 40 (matching the target file) instead of 40.83333...
$errDetailColumnWidth = 40 - (5/6)

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E4").Value = $readyStatus
$ws.Range("F4").Value = $readyStatus
$ws.Range("G4").Value = "2016-09-02 04:27:50"

$ws.Range("E5").Value = $readyStatus
$ws.Range("F5").Value = $readyStatus
$ws.Range("G5").Value = "2016-09-02 04:27:50"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C4").Value = $readyStatus
$ws.Range("H4").Value = "2016-09-02 04:27:43"
$ws.Range("P4").Value = $errBec

$ws.Range("C5").Value = $readyStatus
$ws.Range("H5").Value = "2016-09-02 04:27:43"
$ws.Range("P5").Value = $errCce

$ws.Columns.Item(16).ColumnWidth = $errDetailColumnWidth

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C4").Value = $readyStatus
$ws.Range("H4").Value = "2016-09-02 04:27:50"
$ws.Range("P4").Value = $errBec

$ws.Range("C5").Value = $readyStatus
$ws.Range("H5").Value = "2016-09-02 04:27:50"
$ws.Range("P5").Value = $errCce

$ws.Columns.Item(16).ColumnWidth = $errDetailColumnWidth
